$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 293. This shifts all existing data rows
# 293-392 down to 295-394, matching the new dimension A1:R394.
$ws.Rows.Item(293).EntireRow.Insert()
$ws.Rows.Item(293).EntireRow.Insert()

# New row 293: weekly "Primera" quality price entry for 2021-09-29.
$ws.Range("A293").Value = 8
$ws.Range("B293").Value = "Terminal La Palmera de La Serena"
$ws.Range("C293").Value = "Coquimbo"
$ws.Range("D293").Value = "09/29/2021"
$ws.Range("E293").Value = 4
$ws.Range("F293").Value = 100112023
$ws.Range("G293").Value = "Brócoli"
$ws.Range("H293").Value = "Sin especificar"
$ws.Range("I293").Value = "Primera"
$ws.Range("J293").Value = 3000
$ws.Range("K293").Value = 600
$ws.Range("L293").Value = 700
$ws.Range("M293").Value = 650
$ws.Range("N293").Value = "$/unidad"
$ws.Range("O293").Value = "Provincia del Elquí"
$ws.Range("P293").Value = 650
$ws.Range("Q293").Value = 1
$ws.Range("R293").Value = "Hortaliza"

# New row 294: weekly "Segunda" quality price entry for 2021-09-29.
$ws.Range("A294").Value = 8
$ws.Range("B294").Value = "Terminal La Palmera de La Serena"
$ws.Range("C294").Value = "Coquimbo"
$ws.Range("D294").Value = "09/29/2021"
$ws.Range("E294").Value = 4
$ws.Range("F294").Value = 100112023
$ws.Range("G294").Value = "Brócoli"
$ws.Range("H294").Value = "Sin especificar"
$ws.Range("I294").Value = "Segunda"
$ws.Range("J294").Value = 1600
$ws.Range("K294").Value = 500
$ws.Range("L294").Value = 550
$ws.Range("M294").Value = 525
$ws.Range("N294").Value = "$/unidad"
$ws.Range("O294").Value = "Provincia del Elquí"
$ws.Range("P294").Value = 525
$ws.Range("Q294").Value = 1
$ws.Range("R294").Value = "Hortaliza"
